$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1494.5714
$ws.Range("I86").Value = 1245.3334
$ws.Range("K86").Value = 1245.3334
$ws.Range("M86").Value = -122.3334

$ws.Range("H89").Value = 1494.5714
$ws.Range("I89").Value = 1245.3334
$ws.Range("K89").Value = 6226.666999999999
$ws.Range("M89").Value = -610.6669999999995

$ws.Range("H123").Value = 42110
$ws.Range("J123").Value = 42110
$ws.Range("L123").Value = 42110
$ws.Range("N123").Value = -51910

$ws.Range("H137").Value = 1651.1111
$ws.Range("I137").Value = 1264.6666
$ws.Range("K137").Value = 3793.9998
$ws.Range("M137").Value = -1243.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 900
$ws.Range("I2").Value = 900
$ws.Range("K2").Value = 900
$ws.Range("M2").Value = -787

$ws.Range("H32").Value = 3265.818
$ws.Range("I32").Value = 2114.1887
$ws.Range("K32").Value = 2114.1887
$ws.Range("M32").Value = -1827.1887

$ws.Range("H45").Value = 2317.7
$ws.Range("I45").Value = 1116.4
$ws.Range("K45").Value = 1116.4
$ws.Range("M45").Value = -739.4000000000001

$ws.Range("H61").Value = 2068.5454
$ws.Range("I61").Value = 1363.9
$ws.Range("K61").Value = 1363.9
$ws.Range("M61").Value = -1151.9

$ws.Range("H116").Value = 900
$ws.Range("I116").Value = 900
$ws.Range("K116").Value = 900
$ws.Range("M116").Value = 1394

$ws.Range("H122").Value = 1871.0625
$ws.Range("I122").Value = 1871.0625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5613.1875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3163.1875
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 1729.1904
$ws.Range("I132").Value = 1385
$ws.Range("K132").Value = 4155
$ws.Range("M132").Value = -1625

$ws.Range("H136").Value = 2068.5454
$ws.Range("I136").Value = 1363.9
$ws.Range("K136").Value = 4091.7
$ws.Range("M136").Value = -1541.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 900
$ws.Range("I3").Value = 900
$ws.Range("K3").Value = 900
$ws.Range("M3").Value = -786

$ws.Range("H86").Value = 1527.25
$ws.Range("I86").Value = 1463.3636
$ws.Range("J86").Value = 1605.3334
$ws.Range("K86").Value = 1463.3636
$ws.Range("L86").Value = 1605.3334
$ws.Range("M86").Value = -340.3635999999999
$ws.Range("N86").Value = -3851.3334

$ws.Range("H89").Value = 1527.25
$ws.Range("I89").Value = 1463.3636
$ws.Range("J89").Value = 1605.3334
$ws.Range("K89").Value = 7316.817999999999
$ws.Range("L89").Value = 8026.666999999999
$ws.Range("M89").Value = -1700.817999999999
$ws.Range("N89").Value = -19258.667

$ws.Range("H94").Value = 1909.6
$ws.Range("I94").Value = 1909.6
$ws.Range("K94").Value = 1909.6
$ws.Range("M94").Value = -1458.6

$ws.Range("H105").Value = 2505.2104
$ws.Range("I105").Value = 2505.2104
$ws.Range("K105").Value = 2505.2104
$ws.Range("M105").Value = -758.2103999999999

$ws.Range("H132").Value = 69890
$ws.Range("J132").Value = 69890
$ws.Range("L132").Value = 69890
$ws.Range("N132").Value = -80010

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1433.8276
$ws.Range("I122").Value = 1520.3889
$ws.Range("J122").Value = 1292.1818
$ws.Range("K122").Value = 4561.1667
$ws.Range("L122").Value = 3876.5454
$ws.Range("M122").Value = -2111.1667
$ws.Range("N122").Value = -8776.545399999999

$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -34920

$ws.Range("H134").Value = 831.3125
$ws.Range("I134").Value = 871.5714
$ws.Range("J134").Value = 549.5
$ws.Range("K134").Value = 2614.7142
$ws.Range("L134").Value = 1648.5
$ws.Range("M134").Value = -79.71420000000035
$ws.Range("N134").Value = -6718.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5502953.5
$ws.Range("J131").Value = 9003.035
$ws.Range("L131").Value = 27009.105
$ws.Range("N131").Value = -37089.105

$ws.Range("H132").Value = 1696
$ws.Range("J132").Value = 2035
$ws.Range("L132").Value = 18315
$ws.Range("N132").Value = -23375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2360
$ws.Range("J80").Value = 2433.1428
$ws.Range("L80").Value = 2433.1428
$ws.Range("N80").Value = -4429.1428

$ws.Range("H83").Value = 2360
$ws.Range("J83").Value = 2433.1428
$ws.Range("L83").Value = 12165.714
$ws.Range("N83").Value = -22149.714

$ws.Range("H97").Value = 1186.8462
$ws.Range("I97").Value = 597.5
$ws.Range("J97").Value = 2129.8
$ws.Range("K97").Value = 597.5
$ws.Range("L97").Value = 2129.8
$ws.Range("M97").Value = -101.5
$ws.Range("N97").Value = -3121.8

$ws.Range("H102").Value = 3185.3
$ws.Range("I102").Value = 4489.75
$ws.Range("J102").Value = 2315.6667
$ws.Range("K102").Value = 4489.75
$ws.Range("L102").Value = 2315.6667
$ws.Range("M102").Value = -2867.75
$ws.Range("N102").Value = -5559.6667

$ws.Range("H132").Value = 3853.4
$ws.Range("I132").Value = 3163.5652
$ws.Range("J132").Value = 6120
$ws.Range("K132").Value = 9490.695599999999
$ws.Range("L132").Value = 18360
$ws.Range("M132").Value = -6960.695599999999
$ws.Range("N132").Value = -23420

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6027.9287
$ws.Range("I7").Value = 3083.1667
$ws.Range("J7").Value = 8236.5
$ws.Range("K7").Value = 3083.1667
$ws.Range("L7").Value = 8236.5
$ws.Range("M7").Value = -2971.1667
$ws.Range("N7").Value = -8460.5

$ws.Range("H40").Value = 7438.8
$ws.Range("J40").Value = 12398.111
$ws.Range("L40").Value = 12398.111
$ws.Range("N40").Value = -12670.111

$ws.Range("H55").Value = 194.88889
$ws.Range("I55").Value = 161.42857
$ws.Range("J55").Value = 216.18182
$ws.Range("K55").Value = 161.42857
$ws.Range("L55").Value = 216.18182
$ws.Range("M55").Value = 11.57142999999999
$ws.Range("N55").Value = -562.18182

$ws.Range("H82").Value = 4247.375
$ws.Range("I82").Value = 1966.3334
$ws.Range("J82").Value = 5616
$ws.Range("K82").Value = 1966.3334
$ws.Range("L82").Value = 5616
$ws.Range("M82").Value = -1605.3334
$ws.Range("N82").Value = -6338

$ws.Range("H85").Value = 4247.375
$ws.Range("I85").Value = 1966.3334
$ws.Range("J85").Value = 5616
$ws.Range("K85").Value = 1966.3334
$ws.Range("L85").Value = 5616
$ws.Range("M85").Value = -718.3334
$ws.Range("N85").Value = -8112

$ws.Range("H93").Value = 540
$ws.Range("I93").Value = 193.66667
$ws.Range("K93").Value = 193.66667
$ws.Range("M93").Value = 1054.33333

$ws.Range("H122").Value = 8322.833000000001
$ws.Range("I122").Value = 7186.5713
$ws.Range("J122").Value = 9045.909
$ws.Range("K122").Value = 21559.7139
$ws.Range("L122").Value = 27137.727
$ws.Range("M122").Value = -19109.7139
$ws.Range("N122").Value = -32037.727

$ws.Range("H126").Value = 6027.9287
$ws.Range("I126").Value = 3083.1667
$ws.Range("J126").Value = 8236.5
$ws.Range("K126").Value = 9249.500100000001
$ws.Range("L126").Value = 24709.5
$ws.Range("M126").Value = -6779.500100000001
$ws.Range("N126").Value = -29649.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 58000
$ws.Range("J82").Value = 58000
$ws.Range("L82").Value = 58000
$ws.Range("N82").Value = -58766

$ws.Range("H85").Value = 58000
$ws.Range("J85").Value = 58000
$ws.Range("L85").Value = 58000
$ws.Range("N85").Value = -60652

$ws.Range("H126").Value = 5284.5186
$ws.Range("J126").Value = 7199.625
$ws.Range("L126").Value = 21598.875
$ws.Range("N126").Value = -26538.875

$ws.Range("H132").Value = 2767.88
$ws.Range("I132").Value = 2576.611
$ws.Range("K132").Value = 7729.833
$ws.Range("M132").Value = -5199.833

Write-Output "applied 40 hunks"